$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.027.27'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').Value = '1.728.04'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').Value = '310.41'
$ws.Range('E5').Value = '  -5.39%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '0.4845'
$ws.Range('E7').Value = '  +3.73%  '
$ws.Range('D8').Value = '0.3476'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').Value = '43.53'
$ws.Range('E9').Value = '  +3.06%  '
$ws.Range('D10').Value = '0.07235'
$ws.Range('E10').Value = '  -1.60%  '
$ws.Range('D11').Value = '1.052'
$ws.Range('E11').Value = '  -2.60%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '19.93'
$ws.Range('E13').Value = '  -2.71%  '
$ws.Range('D14').Value = '5.877'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '1.726.73'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D16').Value = '6.831'
$ws.Range('E16').Value = '  -4.43%  '
$ws.Range('D17').Value = '87.01'
$ws.Range('E17').Value = '  -5.51%  '
$ws.Range('D18').Value = '0.00001032'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = '0.06399'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '16.61'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').Value = '5.721'
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').Value = '27.094.08'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').Value = '10.98'
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').Value = '2.075'
$ws.Range('E25').Value = '  -3.67%  '
$ws.Range('D26').Value = '154.12'
$ws.Range('E26').Value = '  -4.92%  '
$ws.Range('D27').Value = '19.96'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = '1.925.24'
$ws.Range('E28').Value = '  -1.64%  '
$ws.Range('D29').Value = '2.071'
$ws.Range('E29').Value = '  -4.35%  '
$ws.Range('D30').Value = '120.64'
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('D31').Value = '1.040'
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('D32').Value = '0.09326'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = '3.645'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '5.385'
$ws.Range('E34').Value = '  -2.67%  '
$ws.Range('D35').Value = '0.05910'
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('D36').Value = '0.02182'
$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = '1.429'
$ws.Range('E37').Value = '  +5.35%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '10.98'
$ws.Range('E38').Value = '  -5.54%  '
$ws.Range('D39').Value = '0.2001'
$ws.Range('E39').Value = '  -3.04%  '
$ws.Range('D40').Value = '4.756'
$ws.Range('E40').Value = '  -2.61%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '0.5980'
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('D43').Value = '1.115'
$ws.Range('E43').Value = '  -5.44%  '
$ws.Range('D44').Value = '7.544'
$ws.Range('E44').Value = '  -2.87%  '
$ws.Range('D45').Value = '12.83'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('E46').Value = '  -4.15%  '
$ws.Range('D47').Value = '0.5618'
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('D48').Value = '119.13'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').Value = '1.851'
$ws.Range('E49').Value = '  -3.64%  '
$ws.Range('D50').Value = '1.103'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('E51').Value = '  -2.20%  '
